# bug fix: catching exceptions from worker thread
# Adds a new Todo task "replace Rabbit with smaller and better divided image"
# (Id 51) to the Active sheet, right after the existing row for Id 21,
# and bumps the "Max Id" tracker on the Config sheet to 51.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Active")

# Insert a new row at position 12, shifting existing rows 12-24 down to 13-25.
$ws.Rows.Item(12).Insert()

$ws.Cells.Item(12, 1).Value = 51
$ws.Cells.Item(12, 2).Value = "replace Rabbit with smaller and better divided image"
$ws.Cells.Item(12, 3).Value = "Todo"
$ws.Cells.Item(12, 4).Value = "Task"

# The "Created" column stores dates as plain text strings (e.g. "8/22/2018"),
# so force text formatting before assigning, then clear the formatting
# override so the cell matches the plain (unstyled) look of its neighbors.
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "8/22/2018"
$ws.Cells.Item(12, 5).ClearFormats()

# Update the "Max Id" value on the Config sheet to reflect the new task id.
$wsConfig = $wb.Worksheets.Item("Config")
$wsConfig.Cells.Item(2, 6).Value = 51
